$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.795.73'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.636.78'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''215.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").Value = '''0.0641'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").Value = '''19.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '1.639.55'
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").Value = '1.862.38'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '0.0₃0776'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").Value = '''63.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '25.818.19'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").Value = '''194.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").Value = '''140.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("E27").Value = '  -5.34%  '
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = '''15.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").Value = '''3.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("D36").Value = '''0.899'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").Value = '''0.553'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").Value = '1.114.14'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.801'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''98.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("D46").Value = '''55.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = '''2.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.83%  '
$ws.Range("D48").Value = '''7.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  -4.30%  '
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("E51").Value = '  -0.35%  '
